# Refactored Parser structure. Fixed problems with reading size of
# classes and methods: every "Number of Lines" value that was incorrectly
# recorded as 0 should actually be 1 (a class/method always has at least
# one line). Once all "0" values are gone, that shared string becomes
# unused and the string table shrinks accordingly.
#
# We reuse the text value "1" that is already present elsewhere in the
# workbook (e.g. classToClassRelations!E3) via copy/paste-special so the
# target cells keep being shared-string text cells (not numbers) without
# disturbing their existing cell style.

$wb = $excel.ActiveWorkbook

$wsClass = $wb.Worksheets.Item("classNumberOfLines")
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")
$wsOne = $wb.Worksheets.Item("classToClassRelations")

$wsOne.Range("E3").Copy()

# classNumberOfLines: column B holds "Number of Lines"
$wsClass.Range("B2").PasteSpecial(-4163)

# methodNumberOfLines: column C holds "Number of Lines"
$rows = @(2, 6, 7, 8, 9, 10, 18, 25, 30, 32, 34, 35, 36, 38, 40, 42)
foreach ($r in $rows) {
    $wsOne.Range("E3").Copy()
    $wsMethod.Cells.Item($r, 3).PasteSpecial(-4163)
}
